# Auto-generated: applies scheduled-runner market-price refresh to Leve profit tables
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 228.57143
$ws.Range("J2").Value = 200
$ws.Range("L2").Value = 200
$ws.Range("N2").Value = -426
$ws.Range("H11").Value = 472.7857
$ws.Range("I11").Value = 472.7857
$ws.Range("K11").Value = 472.7857
$ws.Range("M11").Value = -332.7857
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").Value = 5000
$ws.Range("N16").Value = -5460
$ws.Range("H38").Value = 123.22222
$ws.Range("I38").Value = 123.22222
$ws.Range("K38").Value = 369.66666
$ws.Range("M38").Value = 2.333340000000021
$ws.Range("H43").Value = 2484.4
$ws.Range("I43").Value = 2266.4285
$ws.Range("K43").Value = 2266.4285
$ws.Range("M43").Value = -2197.4285
$ws.Range("H45").Value = 199
$ws.Range("I45").Value = 199
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 597
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -405
$ws.Range("H70").Value = 4999.25
$ws.Range("J70").Value = 4999.25
$ws.Range("L70").Value = 14997.75
$ws.Range("N70").Value = -15537.75
$ws.Range("H73").Value = 4999.25
$ws.Range("J73").Value = 4999.25
$ws.Range("L73").Value = 14997.75
$ws.Range("N73").Value = -16869.75
$ws.Range("H76").Value = 8070.091
$ws.Range("I76").Value = 7642.3335
$ws.Range("K76").Value = 7642.3335
$ws.Range("M76").Value = -7327.3335
$ws.Range("H79").Value = 8070.091
$ws.Range("I79").Value = 7642.3335
$ws.Range("K79").Value = 7642.3335
$ws.Range("M79").Value = -6550.3335
$ws.Range("H86").Value = 2699.8
$ws.Range("J86").Value = 3500
$ws.Range("L86").Value = 3500
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 2699.8
$ws.Range("J89").Value = 3500
$ws.Range("L89").Value = 17500
$ws.Range("N89").Value = -28732
$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492
$ws.Range("H106").Value = 2636.6667
$ws.Range("I106").Value = 1954.5
$ws.Range("K106").Value = 1954.5
$ws.Range("M106").Value = -1323.5
$ws.Range("H112").Value = 2052.7368
$ws.Range("J112").Value = 2082.5293
$ws.Range("L112").Value = 6247.5879
$ws.Range("N112").Value = -8463.5879
$ws.Range("H113").Value = 1887.25
$ws.Range("I113").Value = 1887.25
$ws.Range("K113").Value = 1887.25
$ws.Range("M113").Value = 1366.75
$ws.Range("H116").Value = 3899.875
$ws.Range("I116").Value = 3566.5
$ws.Range("K116").Value = 3566.5
$ws.Range("M116").Value = -124.5
$ws.Range("H129").Value = 12075.667
$ws.Range("I129").Value = 1085.25
$ws.Range("J129").Value = 99999
$ws.Range("K129").Value = 3255.75
$ws.Range("L129").Value = 299997
$ws.Range("M129").Value = 1744.25
$ws.Range("N129").Value = -309997
$ws.Range("H131").Value = 115151.664
$ws.Range("I131").Value = 129233.125
$ws.Range("J131").Value = 2500
$ws.Range("K131").Value = 387699.375
$ws.Range("L131").Value = 7500
$ws.Range("M131").Value = -382659.375
$ws.Range("N131").Value = -17580
$ws.Range("H135").Value = 1837.0869
$ws.Range("I135").Value = 1520.7646
$ws.Range("K135").Value = 13686.8814
$ws.Range("M135").Value = -11151.8814
$ws.Range("H137").Value = 1798.6111
$ws.Range("I137").Value = 1525.2727
$ws.Range("J137").Value = 2228.1428
$ws.Range("K137").Value = 4575.8181
$ws.Range("L137").Value = 6684.428400000001
$ws.Range("M137").Value = -2025.8181
$ws.Range("N137").Value = -11784.4284
$ws.Range("H138").Value = 2388.42
$ws.Range("I138").Value = 1587.0741
$ws.Range("J138").Value = 3329.1304
$ws.Range("K138").Value = 4761.2223
$ws.Range("L138").Value = 9987.3912
$ws.Range("M138").Value = 378.7776999999996
$ws.Range("N138").Value = -20267.3912

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1725.3011
$ws.Range("I32").Value = 1804.1519
$ws.Range("K32").Value = 1804.1519
$ws.Range("M32").Value = -1517.1519
$ws.Range("H37").Value = 17800
$ws.Range("J37").Value = 17800
$ws.Range("L37").Value = 17800
$ws.Range("N37").Value = -18346
$ws.Range("H45").Value = 8017.9546
$ws.Range("I45").Value = 6054.8184
$ws.Range("K45").Value = 6054.8184
$ws.Range("M45").Value = -5677.8184
$ws.Range("H74").Value = 2897.0908
$ws.Range("I74").Value = 1919.7858
$ws.Range("K74").Value = 1919.7858
$ws.Range("M74").Value = -1045.7858
$ws.Range("H77").Value = 2897.0908
$ws.Range("I77").Value = 1919.7858
$ws.Range("K77").Value = 9598.929
$ws.Range("M77").Value = -5230.929
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").Value = 0
$ws.Range("H119").Value = 80000
$ws.Range("J119").Value = 80000
$ws.Range("L119").Value = 80000
$ws.Range("N119").Value = -89676
$ws.Range("H122").Value = 2749.7144
$ws.Range("I122").Value = 2636.9375
$ws.Range("J122").Value = 3110.6
$ws.Range("K122").Value = 7910.8125
$ws.Range("L122").Value = 9331.799999999999
$ws.Range("M122").Value = -5460.8125
$ws.Range("N122").Value = -14231.8
$ws.Range("H132").Value = 4526.359
$ws.Range("I132").Value = 3989.3713
$ws.Range("K132").Value = 11968.1139
$ws.Range("M132").Value = -9438.1139
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").ClearContents()
$ws.Range("N8").Value = 0
$ws.Range("H10").Value = 5000
$ws.Range("I10").Value = 5000
$ws.Range("K10").Value = 5000
$ws.Range("M10").Value = -4860
$ws.Range("H64").Value = 267.5
$ws.Range("J64").Value = 281
$ws.Range("L64").Value = 281
$ws.Range("N64").Value = -731
$ws.Range("H67").Value = 267.5
$ws.Range("J67").Value = 281
$ws.Range("L67").Value = 281
$ws.Range("N67").Value = -1841
$ws.Range("H94").Value = 2204.647
$ws.Range("I94").Value = 1931.2667
$ws.Range("J94").Value = 4255
$ws.Range("K94").Value = 1931.2667
$ws.Range("L94").Value = 4255
$ws.Range("M94").Value = -1480.2667
$ws.Range("N94").Value = -5157
$ws.Range("H105").Value = 1628.1111
$ws.Range("I105").Value = 1639.2941
$ws.Range("K105").Value = 1639.2941
$ws.Range("M105").Value = 107.7058999999999
$ws.Range("H134").Value = 2646.3845
$ws.Range("I134").Value = 2791.4707
$ws.Range("K134").Value = 8374.4121
$ws.Range("M134").Value = -5839.4121

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7921.0557
$ws.Range("J58").Value = 11519.2
$ws.Range("L58").Value = 11519.2
$ws.Range("N58").Value = -11925.2
$ws.Range("H62").Value = 5324.75
$ws.Range("I62").Value = 5324.75
$ws.Range("K62").Value = 5324.75
$ws.Range("M62").Value = -4700.75
$ws.Range("H65").Value = 5324.75
$ws.Range("I65").Value = 5324.75
$ws.Range("K65").Value = 26623.75
$ws.Range("M65").Value = -23503.75
$ws.Range("H107").Value = 648.1
$ws.Range("I107").Value = 513.5
$ws.Range("K107").Value = 513.5
$ws.Range("M107").Value = 1406.5
$ws.Range("H122").Value = 4436.4707
$ws.Range("I122").Value = 4633.3076
$ws.Range("K122").Value = 13899.9228
$ws.Range("M122").Value = -11449.9228
$ws.Range("H132").Value = 3222.7
$ws.Range("I132").Value = 2528.375
$ws.Range("K132").Value = 7585.125
$ws.Range("M132").Value = -5055.125
$ws.Range("H134").Value = 6699.5557
$ws.Range("I134").Value = 2549.3333
$ws.Range("K134").Value = 7647.999899999999
$ws.Range("M134").Value = -5112.999899999999
$ws.Range("H136").Value = 7921.0557
$ws.Range("J136").Value = 11519.2
$ws.Range("L136").Value = 34557.60000000001
$ws.Range("N136").Value = -39657.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4632.8823
$ws.Range("I3").Value = 4647.4375
$ws.Range("K3").Value = 13942.3125
$ws.Range("M3").Value = -13830.3125
$ws.Range("H4").Value = 4936037.5
$ws.Range("I4").Value = 1547817.8
$ws.Range("K4").Value = 4643453.4
$ws.Range("M4").Value = -4643341.4
$ws.Range("H5").Value = 2359.6365
$ws.Range("I5").Value = 2066.4443
$ws.Range("K5").Value = 6199.3329
$ws.Range("M5").Value = -6087.3329
$ws.Range("H68").Value = 1961.25
$ws.Range("I68").Value = 863.6667
$ws.Range("J68").Value = 2619.8
$ws.Range("K68").Value = 2591.0001
$ws.Range("L68").Value = 7859.400000000001
$ws.Range("M68").Value = -1780.0001
$ws.Range("N68").Value = -9481.400000000001
$ws.Range("H71").Value = 1961.25
$ws.Range("I71").Value = 863.6667
$ws.Range("J71").Value = 2619.8
$ws.Range("K71").Value = 7773.0003
$ws.Range("L71").Value = 23578.2
$ws.Range("M71").Value = -3717.0003
$ws.Range("N71").Value = -31690.2
$ws.Range("H109").Value = 3049.8
$ws.Range("I109").Value = 3187.25
$ws.Range("K109").Value = 9561.75
$ws.Range("M109").Value = -8521.75
$ws.Range("H131").Value = 4763809
$ws.Range("I131").Value = 749.25
$ws.Range("J131").Value = 7694922.5
$ws.Range("K131").Value = 2247.75
$ws.Range("L131").Value = 23084767.5
$ws.Range("M131").Value = 2792.25
$ws.Range("N131").Value = -23094847.5
$ws.Range("H132").Value = 1426.8
$ws.Range("I132").Value = 1090.25
$ws.Range("J132").Value = 1811.4286
$ws.Range("K132").Value = 9812.25
$ws.Range("L132").Value = 16302.8574
$ws.Range("M132").Value = -7282.25
$ws.Range("N132").Value = -21362.8574
$ws.Range("H133").Value = 1872.75
$ws.Range("I133").Value = 1872.75
$ws.Range("K133").Value = 5618.25
$ws.Range("M133").Value = -558.25
$ws.Range("H135").Value = 2359.6365
$ws.Range("I135").Value = 2066.4443
$ws.Range("K135").Value = 18597.9987
$ws.Range("M135").Value = -16062.9987
$ws.Range("H136").Value = 364
$ws.Range("I136").Value = 364
$ws.Range("K136").Value = 1092
$ws.Range("M136").Value = 4008
$ws.Range("H138").Value = 777
$ws.Range("I138").Value = 777
$ws.Range("K138").Value = 2331
$ws.Range("M138").Value = 2809

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2080000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 2080000
$ws.Range("K3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value = 2080000
$ws.Range("N3").Value = -2080232
$ws.Range("H70").Value = 6469.8
$ws.Range("I70").Value = 2619.8
$ws.Range("K70").Value = 2619.8
$ws.Range("M70").Value = -2349.8
$ws.Range("H73").Value = 6469.8
$ws.Range("I73").Value = 2619.8
$ws.Range("K73").Value = 2619.8
$ws.Range("M73").Value = -1683.8
$ws.Range("H80").Value = 2827.6333
$ws.Range("I80").Value = 2502.2
$ws.Range("J80").Value = 3153.0667
$ws.Range("K80").Value = 2502.2
$ws.Range("L80").Value = 3153.0667
$ws.Range("M80").Value = -1504.2
$ws.Range("N80").Value = -5149.066699999999
$ws.Range("H83").Value = 2827.6333
$ws.Range("I83").Value = 2502.2
$ws.Range("J83").Value = 3153.0667
$ws.Range("K83").Value = 12511
$ws.Range("L83").Value = 15765.3335
$ws.Range("M83").Value = -7519
$ws.Range("N83").Value = -25749.3335
$ws.Range("H97").Value = 765.5
$ws.Range("I97").Value = 672
$ws.Range("J97").Value = 1233
$ws.Range("K97").Value = 672
$ws.Range("L97").Value = 1233
$ws.Range("M97").Value = -176
$ws.Range("N97").Value = -2225
$ws.Range("H102").Value = 3967.7693
$ws.Range("I102").Value = 3168.2
$ws.Range("K102").Value = 3168.2
$ws.Range("M102").Value = -1546.2
$ws.Range("H122").Value = 1250
$ws.Range("I122").Value = 1250
$ws.Range("K122").Value = 3750
$ws.Range("M122").Value = -1300
$ws.Range("H126").Value = 2611.0952
$ws.Range("I126").Value = 1783.4
$ws.Range("J126").Value = 3363.5454
$ws.Range("K126").Value = 5350.200000000001
$ws.Range("L126").Value = 10090.6362
$ws.Range("M126").Value = -2880.200000000001
$ws.Range("N126").Value = -15030.6362
$ws.Range("H132").Value = 3255.3572
$ws.Range("I132").Value = 2447.125
$ws.Range("K132").Value = 7341.375
$ws.Range("M132").Value = -4811.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4959
$ws.Range("I40").Value = 4959
$ws.Range("K40").Value = 4959
$ws.Range("M40").Value = -4823
$ws.Range("H68").Value = 13040.667
$ws.Range("I68").Value = 13695.477
$ws.Range("J68").Value = 10748.833
$ws.Range("K68").Value = 13695.477
$ws.Range("L68").Value = 10748.833
$ws.Range("M68").Value = -12946.477
$ws.Range("N68").Value = -12246.833
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H71").Value = 13040.667
$ws.Range("I71").Value = 13695.477
$ws.Range("J71").Value = 10748.833
$ws.Range("K71").Value = 68477.38500000001
$ws.Range("L71").Value = 53744.165
$ws.Range("M71").Value = -64733.38500000001
$ws.Range("N71").Value = -61232.165
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H93").Value = 16901.8
$ws.Range("J93").Value = 27558.334
$ws.Range("L93").Value = 27558.334
$ws.Range("N93").Value = -30054.334
$ws.Range("H122").Value = 2776.2
$ws.Range("I122").Value = 2397.7144
$ws.Range("J122").Value = 3107.375
$ws.Range("K122").Value = 7193.1432
$ws.Range("L122").Value = 9322.125
$ws.Range("M122").Value = -4743.1432
$ws.Range("N122").Value = -14222.125
$ws.Range("H132").Value = 6354.533
$ws.Range("I132").Value = 6354.533
$ws.Range("K132").Value = 19063.599
$ws.Range("M132").Value = -16533.599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 15000
$ws.Range("J49").Value = 15000
$ws.Range("L49").Value = 15000
$ws.Range("N49").Value = -15460
$ws.Range("H107").Value = 1142.24
$ws.Range("I107").Value = 680
$ws.Range("K107").Value = 2040
$ws.Range("M107").Value = -120
$ws.Range("H132").Value = 3282.8604
$ws.Range("I132").Value = 2406.353
$ws.Range("K132").Value = 7219.059
$ws.Range("M132").Value = -4689.059
$ws.Range("H136").Value = 5541.6665
$ws.Range("I136").Value = 4681.8184
$ws.Range("K136").Value = 14045.4552
$ws.Range("M136").Value = -11495.4552

Write-Host "Applied market-price refresh across all sheets."